$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("optimization_parameters")

# Remove the duplicated "value" header cells in row 1 (C1:F1), leaving only A1/B1.
$ws.Range("C1:F1").ClearContents() | Out-Null

# Remove the obsolete "Deletion" parameter row (old row 17); this shifts the
# "simulation_timepoints" row (old row 18) up to become the new row 17, and
# drops the now-unused "Deletion" shared string from the workbook.
$ws.Rows.Item(17).Delete() | Out-Null

# Update the selected/active cell shown when the sheet is opened.
$ws.Range("G4").Select() | Out-Null
